# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) on Sheet1 held the literal text "6-8-2012-13"
# (the source filename, not a real date) for every data row; correct it
# to the actual game date "2013-06-08", keeping the value as plain text
# (not an Excel date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Text -eq "6-8-2012-13") {
        # Assigning the literal string directly would make Excel
        # reinterpret the "YYYY-MM-DD"-looking text as a date serial
        # number. Instead, build it as a text-formula result (guaranteed
        # string) and paste that back over itself as a value, which
        # keeps the cell a plain text cell without touching its
        # NumberFormat/style.
        $cell.Formula = '="2013-06-08"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}
